$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "MAE" column before the existing "Tipo" column, pushing
# "Tipo" from D to E.
$ws.Range("D1").EntireColumn.Insert()

# Header for the new column, styled like the rest of the header row.
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Updated metric values for the single AdaBoostRegressor row.
$ws.Range("B2").Value = 0.1524307648906948
$ws.Range("C2").Value = 0.998428928608553
$ws.Range("D2").Value = 0.3316364171675909
